# Fix: new A11. 五岁以下死亡监测(季度) DeathMonitorUnder5QuarterlyArea
# Row 11 on sheet "列表" holds three label/value pairs for 填报单位(reporting
# unit) / 填报人(reporter) / 填报日期(report date). The "填报单位" caption
# cell was dropped and the merged value field shifted one column to the
# left so that it starts right under the row-10 grouping; the "填报人"
# caption/value pair shifted left by one column the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- "填报单位" block: was B11 (caption) + merged C11:E11 (value) -------
$fuwujgmc = $ws.Range("C11").Value2
$ws.Range("C11:E11").UnMerge()
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("B11").Value = $fuwujgmc
$ws.Range("B11:D11").Merge()

# --- "填报人" block: was J11 (caption) + merged K11:M11 (value) --------
$tianbaoren = $ws.Range("J11").Value2
$dataUserName = $ws.Range("K11").Value2
$ws.Range("K11:M11").UnMerge()
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("I11").Value = $tianbaoren
$ws.Range("J11").Value = $dataUserName
$ws.Range("J11:L11").Merge()
